$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("factors")

# Update the "Done?" column (D) values to "Yes" for the rows that were
# previously marked "Possible" / "Done", and fill in the previously blank
# D13 cell with "Yes" as well.
$ws.Range("D6").Value  = "Yes"
$ws.Range("D13").Value = "Yes"
$ws.Range("D19").Value = "Yes"
$ws.Range("D24").Value = "Yes"
$ws.Range("D26").Value = "Yes"
$ws.Range("D27").Value = "Yes"
$ws.Range("D28").Value = "Yes"

# Scroll the view down a bit and move the selection, matching the saved
# sheet view state (topLeftCell D5, selection E23).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E23").Select()
